# Auto-generated edit script: updates leve-profit calculation cells
# across ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets per the commit diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(62, 8).Value = 3800  # H62: 5225.625 -> 3800
$ws.Cells.Item(62, 9).Value = 1933.3334  # I62: 1951.25 -> 1933.3334
$ws.Cells.Item(62, 10).Value = 15000  # J62: 8500 -> 15000
$ws.Cells.Item(62, 11).Value = 1933.3334  # K62: 1951.25 -> 1933.3334
$ws.Cells.Item(62, 12).Value = 15000  # L62: 8500 -> 15000
$ws.Cells.Item(62, 13).Value = -1309.3334  # M62: -1327.25 -> -1309.3334
$ws.Cells.Item(62, 14).Value = -16248  # N62: -9748 -> -16248
$ws.Cells.Item(65, 8).Value = 3800  # H65: 5225.625 -> 3800
$ws.Cells.Item(65, 9).Value = 1933.3334  # I65: 1951.25 -> 1933.3334
$ws.Cells.Item(65, 10).Value = 15000  # J65: 8500 -> 15000
$ws.Cells.Item(65, 11).Value = 9666.666999999999  # K65: 9756.25 -> 9666.666999999999
$ws.Cells.Item(65, 12).Value = 75000  # L65: 42500 -> 75000
$ws.Cells.Item(65, 13).Value = -6546.666999999999  # M65: -6636.25 -> -6546.666999999999
$ws.Cells.Item(65, 14).Value = -81240  # N65: -48740 -> -81240
$ws.Cells.Item(100, 8).Value = 20001326  # H100: 20001420 -> 20001326
$ws.Cells.Item(100, 9).Value = 22223472  # I100: 20001420 -> 22223472
$ws.Cells.Item(100, 10).Value = 2000  # J100: 0 -> 2000
$ws.Cells.Item(100, 11).Value = 22223472  # K100: 20001420 -> 22223472
$ws.Cells.Item(100, 12).Value = 2000  # L100: 0 -> 2000
$ws.Cells.Item(100, 13).Value = -22222931  # M100: -20000879 -> -22222931
$ws.Cells.Item(100, 14).Value = -3082  # N100: None -> -3082
$ws.Cells.Item(112, 8).Value = 10001554  # H112: 10418229 -> 10001554
$ws.Cells.Item(112, 10).Value = 1578.1837  # J112: 1588.1277 -> 1578.1837
$ws.Cells.Item(112, 12).Value = 4734.551100000001  # L112: 4764.3831 -> 4734.551100000001
$ws.Cells.Item(112, 14).Value = -6950.551100000001  # N112: -6980.3831 -> -6950.551100000001
$ws.Cells.Item(132, 8).Value = 22313350  # H132: 21828570 -> 22313350
$ws.Cells.Item(132, 9).Value = 23350622  # I132: 23906660 -> 23350622
$ws.Cells.Item(132, 10).Value = 12000  # J132: 8625 -> 12000
$ws.Cells.Item(132, 11).Value = 70051866  # K132: 71719980 -> 70051866
$ws.Cells.Item(132, 12).Value = 36000  # L132: 25875 -> 36000
$ws.Cells.Item(132, 13).Value = -70049336  # M132: -71717450 -> -70049336
$ws.Cells.Item(132, 14).Value = -41060  # N132: -30935 -> -41060
$ws.Cells.Item(137, 8).Value = 2893.0625  # H137: 3160.2144 -> 2893.0625
$ws.Cells.Item(137, 9).Value = 2170.4285  # I137: 2338.1667 -> 2170.4285
$ws.Cells.Item(137, 10).Value = 4272.636  # J137: 4639.9 -> 4272.636
$ws.Cells.Item(137, 11).Value = 6511.2855  # K137: 7014.500100000001 -> 6511.2855
$ws.Cells.Item(137, 12).Value = 12817.908  # L137: 13919.7 -> 12817.908
$ws.Cells.Item(137, 13).Value = -3961.2855  # M137: -4464.500100000001 -> -3961.2855
$ws.Cells.Item(137, 14).Value = -17917.908  # N137: -19019.7 -> -17917.908
$ws.Cells.Item(138, 8).Value = 2257.42  # H138: 2332.53 -> 2257.42
$ws.Cells.Item(138, 9).Value = 863.37036  # I138: 958.3103599999999 -> 863.37036
$ws.Cells.Item(138, 10).Value = 2773.0273  # J138: 2893.831 -> 2773.0273
$ws.Cells.Item(138, 11).Value = 2590.11108  # K138: 2874.93108 -> 2590.11108
$ws.Cells.Item(138, 12).Value = 8319.081900000001  # L138: 8681.493 -> 8319.081900000001
$ws.Cells.Item(138, 13).Value = 2549.88892  # M138: 2265.06892 -> 2549.88892
$ws.Cells.Item(138, 14).Value = -18599.0819  # N138: -18961.493 -> -18599.0819

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 7695.803  # H32: 7284.3286 -> 7695.803
$ws.Cells.Item(32, 9).Value = 4813.59  # I32: 4411.8604 -> 4813.59
$ws.Cells.Item(32, 11).Value = 4813.59  # K32: 4411.8604 -> 4813.59
$ws.Cells.Item(32, 13).Value = -4526.59  # M32: -4124.8604 -> -4526.59
$ws.Cells.Item(45, 8).Value = 1388  # H45: 1267 -> 1388
$ws.Cells.Item(45, 9).Value = 1204  # I45: 1124 -> 1204
$ws.Cells.Item(45, 10).Value = 1480  # J45: 1410 -> 1480
$ws.Cells.Item(45, 11).Value = 1204  # K45: 1124 -> 1204
$ws.Cells.Item(45, 12).Value = 1480  # L45: 1410 -> 1480
$ws.Cells.Item(45, 13).Value = -827  # M45: -747 -> -827
$ws.Cells.Item(45, 14).Value = -2234  # N45: -2164 -> -2234

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(29, 8).Value = 6993  # H29: 3999 -> 6993
$ws.Cells.Item(29, 9).Value = 2989.5  # I29: 1998.75 -> 2989.5
$ws.Cells.Item(29, 10).Value = 15000  # J29: 12000 -> 15000
$ws.Cells.Item(29, 11).Value = 2989.5  # K29: 1998.75 -> 2989.5
$ws.Cells.Item(29, 12).Value = 15000  # L29: 12000 -> 15000
$ws.Cells.Item(29, 13).Value = -2700.5  # M29: -1709.75 -> -2700.5
$ws.Cells.Item(29, 14).Value = -15578  # N29: -12578 -> -15578

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 29416542  # H31: 17245128 -> 29416542
$ws.Cells.Item(31, 9).Value = 1353.3334  # I31: 1570.9231 -> 1353.3334
$ws.Cells.Item(31, 10).Value = 45461190  # J31: 31255518 -> 45461190
$ws.Cells.Item(31, 11).Value = 1353.3334  # K31: 1570.9231 -> 1353.3334
$ws.Cells.Item(31, 12).Value = 45461190  # L31: 31255518 -> 45461190
$ws.Cells.Item(31, 13).Value = -1058.3334  # M31: -1275.9231 -> -1058.3334
$ws.Cells.Item(31, 14).Value = -45461780  # N31: -31256108 -> -45461780
$ws.Cells.Item(34, 8).Value = 29416542  # H34: 17245128 -> 29416542
$ws.Cells.Item(34, 9).Value = 1353.3334  # I34: 1570.9231 -> 1353.3334
$ws.Cells.Item(34, 10).Value = 45461190  # J34: 31255518 -> 45461190
$ws.Cells.Item(34, 11).Value = 1353.3334  # K34: 1570.9231 -> 1353.3334
$ws.Cells.Item(34, 12).Value = 45461190  # L34: 31255518 -> 45461190
$ws.Cells.Item(34, 13).Value = -1151.3334  # M34: -1368.9231 -> -1151.3334
$ws.Cells.Item(34, 14).Value = -45461594  # N34: -31255922 -> -45461594
$ws.Cells.Item(58, 8).Value = 1767.9375  # H58: 1792.359 -> 1767.9375
$ws.Cells.Item(58, 9).Value = 1626.6451  # I58: 1636.3934 -> 1626.6451
$ws.Cells.Item(58, 10).Value = 2254.611  # J58: 2352 -> 2254.611
$ws.Cells.Item(58, 11).Value = 1626.6451  # K58: 1636.3934 -> 1626.6451
$ws.Cells.Item(58, 12).Value = 2254.611  # L58: 2352 -> 2254.611
$ws.Cells.Item(58, 13).Value = -1423.6451  # M58: -1433.3934 -> -1423.6451
$ws.Cells.Item(58, 14).Value = -2660.611  # N58: -2758 -> -2660.611
$ws.Cells.Item(86, 8).Value = 3000  # H86: 3666.6667 -> 3000
$ws.Cells.Item(86, 10).Value = 3000  # J86: 3666.6667 -> 3000
$ws.Cells.Item(86, 12).Value = 3000  # L86: 3666.6667 -> 3000
$ws.Cells.Item(86, 14).Value = -5246  # N86: -5912.6667 -> -5246
$ws.Cells.Item(89, 8).Value = 3000  # H89: 3666.6667 -> 3000
$ws.Cells.Item(89, 10).Value = 3000  # J89: 3666.6667 -> 3000
$ws.Cells.Item(89, 12).Value = 15000  # L89: 18333.3335 -> 15000
$ws.Cells.Item(89, 14).Value = -26232  # N89: -29565.3335 -> -26232
$ws.Cells.Item(122, 8).Value = 1893.9565  # H122: 2019.85 -> 1893.9565
$ws.Cells.Item(122, 9).Value = 1138.0588  # I122: 1180.1875 -> 1138.0588
$ws.Cells.Item(122, 10).Value = 4035.6667  # J122: 5378.5 -> 4035.6667
$ws.Cells.Item(122, 11).Value = 3414.1764  # K122: 3540.5625 -> 3414.1764
$ws.Cells.Item(122, 12).Value = 12107.0001  # L122: 16135.5 -> 12107.0001
$ws.Cells.Item(122, 13).Value = -964.1764000000003  # M122: -1090.5625 -> -964.1764000000003
$ws.Cells.Item(122, 14).Value = -17007.0001  # N122: -21035.5 -> -17007.0001
$ws.Cells.Item(132, 8).Value = 3883.0833  # H132: 3123.6287 -> 3883.0833
$ws.Cells.Item(132, 9).Value = 1812.8889  # I132: 1414.5625 -> 1812.8889
$ws.Cells.Item(132, 10).Value = 5125.2  # J132: 4562.8423 -> 5125.2
$ws.Cells.Item(132, 11).Value = 5438.6667  # K132: 4243.6875 -> 5438.6667
$ws.Cells.Item(132, 12).Value = 15375.6  # L132: 13688.5269 -> 15375.6
$ws.Cells.Item(132, 13).Value = -2908.6667  # M132: -1713.6875 -> -2908.6667
$ws.Cells.Item(132, 14).Value = -20435.6  # N132: -18748.5269 -> -20435.6
$ws.Cells.Item(136, 8).Value = 1767.9375  # H136: 1792.359 -> 1767.9375
$ws.Cells.Item(136, 9).Value = 1626.6451  # I136: 1636.3934 -> 1626.6451
$ws.Cells.Item(136, 10).Value = 2254.611  # J136: 2352 -> 2254.611
$ws.Cells.Item(136, 11).Value = 4879.9353  # K136: 4909.1802 -> 4879.9353
$ws.Cells.Item(136, 12).Value = 6763.833  # L136: 7056 -> 6763.833
$ws.Cells.Item(136, 13).Value = -2329.9353  # M136: -2359.1802 -> -2329.9353
$ws.Cells.Item(136, 14).Value = -11863.833  # N136: -12156 -> -11863.833

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(86, 8).Value = 10959.667  # H86: 6746.923 -> 10959.667
$ws.Cells.Item(86, 9).Value = 25125  # I86: 7274.75 -> 25125
$ws.Cells.Item(86, 10).Value = 3877  # J86: 5902.4 -> 3877
$ws.Cells.Item(86, 11).Value = 75375  # K86: 21824.25 -> 75375
$ws.Cells.Item(86, 12).Value = 11631  # L86: 17707.2 -> 11631
$ws.Cells.Item(86, 13).Value = -74189  # M86: -20638.25 -> -74189
$ws.Cells.Item(86, 14).Value = -14003  # N86: -20079.2 -> -14003
$ws.Cells.Item(89, 8).Value = 10959.667  # H89: 6746.923 -> 10959.667
$ws.Cells.Item(89, 9).Value = 25125  # I89: 7274.75 -> 25125
$ws.Cells.Item(89, 10).Value = 3877  # J89: 5902.4 -> 3877
$ws.Cells.Item(89, 11).Value = 226125  # K89: 65472.75 -> 226125
$ws.Cells.Item(89, 12).Value = 34893  # L89: 53121.6 -> 34893
$ws.Cells.Item(89, 13).Value = -220197  # M89: -59544.75 -> -220197
$ws.Cells.Item(89, 14).Value = -46749  # N89: -64977.6 -> -46749
$ws.Cells.Item(113, 8).Value = 581.0857  # H113: 583.4706 -> 581.0857
$ws.Cells.Item(113, 9).Value = 561.1539  # I113: 563.6 -> 561.1539
$ws.Cells.Item(113, 11).Value = 1683.4617  # K113: 1690.8 -> 1683.4617
$ws.Cells.Item(113, 13).Value = 486.5382999999999  # M113: 479.1999999999998 -> 486.5382999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(48, 8).Value = 34999.5  # H48: 35000 -> 34999.5
$ws.Cells.Item(48, 10).Value = 34999.5  # J48: 35000 -> 34999.5
$ws.Cells.Item(48, 12).Value = 34999.5  # L48: 35000 -> 34999.5
$ws.Cells.Item(48, 14).Value = -35969.5  # N48: -35970 -> -35969.5
$ws.Cells.Item(107, 8).Value = 850.3  # H107: 915 -> 850.3
$ws.Cells.Item(107, 9).Value = 581.6667  # I107: 623 -> 581.6667
$ws.Cells.Item(107, 10).Value = 1253.25  # J107: 2083 -> 1253.25
$ws.Cells.Item(107, 11).Value = 581.6667  # K107: 623 -> 581.6667
$ws.Cells.Item(107, 12).Value = 1253.25  # L107: 2083 -> 1253.25
$ws.Cells.Item(107, 13).Value = 1338.3333  # M107: 1297 -> 1338.3333
$ws.Cells.Item(107, 14).Value = -5093.25  # N107: -5923 -> -5093.25
$ws.Cells.Item(122, 8).Value = 2984.5833  # H122: 2778 -> 2984.5833
$ws.Cells.Item(122, 9).Value = 2225.875  # I122: 1827.8182 -> 2225.875
$ws.Cells.Item(122, 10).Value = 4502  # J122: 8004 -> 4502
$ws.Cells.Item(122, 11).Value = 6677.625  # K122: 5483.4546 -> 6677.625
$ws.Cells.Item(122, 12).Value = 13506  # L122: 24012 -> 13506
$ws.Cells.Item(122, 13).Value = -4227.625  # M122: -3033.4546 -> -4227.625
$ws.Cells.Item(122, 14).Value = -18406  # N122: -28912 -> -18406
$ws.Cells.Item(126, 8).Value = 2942.89  # H126: 4257.675 -> 2942.89
$ws.Cells.Item(126, 9).Value = 2928.5613  # I126: 2902.7778 -> 2928.5613
$ws.Cells.Item(126, 10).Value = 3645  # J126: 5366.227 -> 3645
$ws.Cells.Item(126, 11).Value = 8785.6839  # K126: 8708.3334 -> 8785.6839
$ws.Cells.Item(126, 12).Value = 10935  # L126: 16098.681 -> 10935
$ws.Cells.Item(126, 13).Value = -6315.6839  # M126: -6238.3334 -> -6315.6839
$ws.Cells.Item(126, 14).Value = -15875  # N126: -21038.681 -> -15875
$ws.Cells.Item(137, 8).Value = 63301.6  # H137: 69234.5 -> 63301.6
$ws.Cells.Item(137, 10).Value = 63301.6  # J137: 69234.5 -> 63301.6
$ws.Cells.Item(137, 12).Value = 63301.6  # L137: 69234.5 -> 63301.6
$ws.Cells.Item(137, 14).Value = -73501.60000000001  # N137: -79434.5 -> -73501.60000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(136, 8).Value = 3979.7036  # H136: 3998.1538 -> 3979.7036
$ws.Cells.Item(136, 10).Value = 7310.5  # J136: 7733.8887 -> 7310.5
$ws.Cells.Item(136, 12).Value = 21931.5  # L136: 23201.6661 -> 21931.5
$ws.Cells.Item(136, 14).Value = -27031.5  # N136: -28301.6661 -> -27031.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 13891375  # H132: 19611004 -> 13891375
$ws.Cells.Item(132, 9).Value = 972.1111  # I132: 1000 -> 972.1111
$ws.Cells.Item(132, 10).Value = 22225618  # J132: 20836628 -> 22225618
$ws.Cells.Item(132, 11).Value = 2916.3333  # K132: 3000 -> 2916.3333
$ws.Cells.Item(132, 12).Value = 66676854  # L132: 62509884 -> 66676854
$ws.Cells.Item(132, 13).Value = -386.3332999999998  # M132: -470 -> -386.3332999999998
$ws.Cells.Item(132, 14).Value = -66681914  # N132: -62514944 -> -66681914
